$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Update Disease Ontology (row 3) version from v2024-12-18 -> v2025-02-03
$ws.Range("E3").Value = "v2025-02-03"

# Update Experimental Factor Ontology (row 4) version from v3.73.0 -> v3.74.0
$ws.Range("E4").Value = "v3.74.0"

# Move the active selection from E3 to E4
$ws.Range("E4").Select()
